# Regenerate orders with updated distance/size labels.
# Mapping (applied as substring replacement, order matters to avoid collisions):
#   D64 -> D69
#   D51 -> D55
#   D80 -> D86
#   S30 -> S31
# (S20 / S25 are left unchanged)

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$used = $ws.UsedRange
$rowCount = $used.Rows.Count
$colCount = $used.Columns.Count

for ($r = 1; $r -le $rowCount; $r++) {
    for ($c = 1; $c -le $colCount; $c++) {
        $cell = $ws.Cells.Item($r, $c)
        $val = $cell.Value2

        if ($val -is [string]) {
            $newVal = $val
            $newVal = $newVal.Replace("D64", "D69")
            $newVal = $newVal.Replace("D51", "D55")
            $newVal = $newVal.Replace("D80", "D86")
            $newVal = $newVal.Replace("S30", "S31")

            if ($newVal -ne $val) {
                $cell.Value2 = $newVal
            }
        }
    }
}
